$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "synthetic"
$ws.Range("C7").Value = "original"
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = "C:\Users\franz\Documents\work\projects\arp\data\synthetic_data\synthetic_data_original_untextured_unclipped_vtp_paraview"
$ws.Range("F7").Value = ".vtp"
$ws.Range("G7").Value = 200
$ws.Range("H7").Value = "M"
$ws.Range("I7").Value = $false
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = "automatic"
$ws.Range("M7").Value = $true
$ws.Range("N7").Value = $false
$ws.Range("O7").Value = $true

$ws.Range("N7").Select()
